# Added EV to DE00. Added some batteries.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 previously only held a lone, empty-valued B6 cell (s="3").
# Fill it in with the new DE00 / Distributed Energy / 2030 / 250000 entry.
$ws.Range("A6").Value = "DE00"
$ws.Range("B6").Value = "Distributed Energy"
$ws.Range("C6").Value = 2030
$ws.Range("D6").Value = 250000

# Add the matching 2040 entry as a brand new row 7.
$ws.Range("A7").Value = "DE00"
$ws.Range("B7").Value = "Distributed Energy"
$ws.Range("C7").Value = 2040
$ws.Range("D7").Value = 500000

# Widen column B (now populated with "Distributed Energy" down through row 7).
$ws.Columns.Item(2).ColumnWidth = 15.67

# Leave the cursor where the author left it after entering the new rows.
$ws.Range("A9").Select()
